$wb = $excel.ActiveWorkbook

# --- Sheet "RCA": insert a new TSL (timeslice) row for the SEASON level ---
$ws = $wb.Worksheets.Item("RCA")

# Insert a new row above row 9. Excel copies the formatting of the row
# above (row 8) down into the freshly inserted row.
$ws.Rows("9:9").Insert()

# Row 8 already described the NRG/RSDHET DAYNITE-less entry - tag it with
# the DAYNITE timeslice level.
$ws.Range("S8").Value = "DAYNITE"

# The newly inserted row 9 becomes a copy of that commodity entry, at the
# SEASON timeslice level.
$ws.Range("M9").Value = "NRG"
$ws.Range("O9").Value = "RSDHET"
$ws.Range("P9").Value = "Residential heat from district heating network copy"
$ws.Range("Q9").Value = "PJ"
$ws.Range("S9").Value = "SEASON"

# Make "RCA" the active sheet/tab, with the same selection/scroll state
# the author left the workbook in.
$ws.Activate() | Out-Null
$ws.Range("S10").Select() | Out-Null
